# fix global slide_num error & make available Exporting slides in for loop
#
# Replaces the placeholder lyric text on slide 1 and slide 2 with the
# real Korean lyrics, run by run, so the existing <a:r>/<a:br/>/<a:r>
# paragraph structure (and each run's formatting) is preserved exactly.

$p = $ppt.ActivePresentation

function Set-LyricRuns {
    param($Slide, $FirstLine, $SecondLine)

    $shape = $Slide.Shapes.Item(1)
    $frame = $shape.TextFrame

    # Remember the shape's original (pre-autofit-relayout) height so we
    # can restore it after editing the text - editing via TextRange
    # triggers the textbox's spAutoFit relayout (height shrinks to fit
    # the new text), which we don't want to leak into the saved geometry.
    $origHeight = $shape.Height

    $textRange = $frame.TextRange
    $fullText  = $textRange.Text

    # The paragraph is "<run1><br/><run2>" - as plain text that is
    # run1 + Chr(11) + run2 (Chr(11) = PowerPoint's vertical-tab line break).
    $breakPos = $fullText.IndexOf([char]11)
    $firstLen = $breakPos

    # Replace the first run's characters in place.
    $run1 = $textRange.Characters(1, $firstLen)
    $run1.Text = $FirstLine

    # Recompute where the second run starts now that the first run's
    # length may have changed (Korean replacement text differs in length
    # from the original placeholder).
    $afterFirst   = $textRange.Text
    $newBreakPos  = $afterFirst.IndexOf([char]11)
    $secondStart  = $newBreakPos + 2
    $secondLen    = $afterFirst.Length - $newBreakPos - 1

    $run2 = $textRange.Characters($secondStart, $secondLen)
    $run2.Text = $SecondLine

    # Restore the shape's original height (undo the autofit relayout).
    $shape.Height = $origHeight
}

Set-LyricRuns $p.Slides.Item(1) "너에게 난" "한여름 노을처럼"
Set-LyricRuns $p.Slides.Item(2) "한편의 아름다운" "추억이 되고"
